$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-61) holds the "Förändrad" date, stored as serial 45205 (2023-10-06).
# Update it to serial 45206 (2023-10-07) for every data row.
for ($row = 2; $row -le 61; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
